$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Preserve the existing cell style while forcing the value to be
    # stored as text (avoids Excel auto-converting numeric-looking
    # strings like '318.70' or '41.582.80' into numbers and losing
    # formatting / trailing zeros / thousand-dot separators).
    $existingStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $existingStyle
}

# Row 2
Set-TextCell $ws.Range('D2') '41.582.80'
$ws.Range('E2').Value = '  -0.22%  '

# Row 3
Set-TextCell $ws.Range('D3') '2.460.58'
$ws.Range('E3').Value = '  -0.41%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
Set-TextCell $ws.Range('D5') '318.70'
$ws.Range('E5').Value = '  +0.55%  '

# Row 6
Set-TextCell $ws.Range('D6') '91.28'
$ws.Range('E6').Value = '  -1.54%  '

# Row 7
Set-TextCell $ws.Range('D7') '0.548'
$ws.Range('E7').Value = '  -0.46%  '

# Row 8
$ws.Range('E8').Value = '  +0.04%  '

# Row 9
Set-TextCell $ws.Range('D9') '0.504'
$ws.Range('E9').Value = '  -2.02%  '

# Row 10
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws.Range('D10') '0.0851'
$ws.Range('E10').Value = '  -4.99%  '

# Row 11
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws.Range('D11') '32.64'
$ws.Range('E11').Value = '  -0.60%  '

# Row 12
$ws.Range('E12').Value = '  -0.97%  '

# Row 13
Set-TextCell $ws.Range('D13') '2.839.37'
$ws.Range('E13').Value = '  -0.48%  '

# Row 14
Set-TextCell $ws.Range('D14') '6.84'
$ws.Range('E14').Value = '  -0.88%  '

# Row 15
Set-TextCell $ws.Range('D15') '15.43'
$ws.Range('E15').Value = '  -1.86%  '

# Row 16
Set-TextCell $ws.Range('D16') '2.490.42'
$ws.Range('E16').Value = '  +0.11%  '

# Row 17
Set-TextCell $ws.Range('D17') '0.785'
$ws.Range('E17').Value = '  +0.09%  '

# Row 18
Set-TextCell $ws.Range('D18') '41.524.36'
$ws.Range('E18').Value = '  -0.37%  '

# Row 19
Set-TextCell $ws.Range('D19') '6.41'
$ws.Range('E19').Value = '  -1.37%  '

# Row 20
Set-TextCell $ws.Range('D20') '0.0₃0937'
$ws.Range('E20').Value = '  -4.16%  '

# Row 21
Set-TextCell $ws.Range('D21') '72.11'
$ws.Range('E21').Value = '  +1.18%  '

# Row 22
Set-TextCell $ws.Range('D22') '11.15'
$ws.Range('E22').Value = '  -2.56%  '

# Row 23
Set-TextCell $ws.Range('D23') '238.24'
$ws.Range('E23').Value = '  -0.57%  '

# Row 24
$ws.Range('E24').Value = '  +0.43%  '

# Row 25
$ws.Range('E25').Value = '  +0.49%  '

# Row 26
$ws.Range('E26').Value = '  +0.07%  '

# Row 27
Set-TextCell $ws.Range('D27') '24.65'
$ws.Range('E27').Value = '  -0.31%  '

# Row 28
Set-TextCell $ws.Range('D28') '2.24'
$ws.Range('E28').Value = '  -1.59%  '

# Row 29
Set-TextCell $ws.Range('D29') '9.66'
$ws.Range('E29').Value = '  -1.61%  '

# Row 30
Set-TextCell $ws.Range('D30') '36.07'
$ws.Range('E30').Value = '  +2.24%  '

# Row 31
Set-TextCell $ws.Range('D31') '156.16'
$ws.Range('E31').Value = '  +0.05%  '

# Row 32
Set-TextCell $ws.Range('D32') '5.40'
$ws.Range('E32').Value = '  -1.89%  '

# Row 33
$ws.Range('E33').Value = '  -0.07%  '

# Row 34
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws.Range('D34') '2.57'
$ws.Range('E34').Value = '  -0.49%  '

# Row 35
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range('D35') '0.0761'
$ws.Range('E35').Value = '  -0.87%  '

# Row 36
Set-TextCell $ws.Range('D36') '16.86'
$ws.Range('E36').Value = '  -3.58%  '

# Row 37
Set-TextCell $ws.Range('D37') '2.89'

# Row 38
$ws.Range('E38').Value = '  +0.17%  '

# Row 39
Set-TextCell $ws.Range('D39') '1.82'
$ws.Range('E39').Value = '  +0.79%  '

# Row 40
$ws.Range('E40').Value = '  -0.45%  '

# Row 41
Set-TextCell $ws.Range('D41') '3.98'
$ws.Range('E41').Value = '  -0.30%  '

# Row 42
Set-TextCell $ws.Range('D42') '2.31'
$ws.Range('E42').Value = '  -7.93%  '

# Row 43
Set-TextCell $ws.Range('D43') '1.999.26'
$ws.Range('E43').Value = '  +1.54%  '

# Row 44
Set-TextCell $ws.Range('D44') '0.0280'
$ws.Range('E44').Value = '  -1.46%  '

# Row 45
Set-TextCell $ws.Range('D45') '18.49'
$ws.Range('E45').Value = '  -1.89%  '

# Row 46
$ws.Range('E46').Value = '  -0.39%  '

# Row 47
Set-TextCell $ws.Range('D47') '9.52'
$ws.Range('E47').Value = '  +4.82%  '

# Row 48
Set-TextCell $ws.Range('D48') '2.724.71'
$ws.Range('E48').Value = '  +0.64%  '

# Row 49
Set-TextCell $ws.Range('D49') '75.79'
$ws.Range('E49').Value = '  +4.29%  '

# Row 50
Set-TextCell $ws.Range('D50') '96.76'
$ws.Range('E50').Value = '  -0.64%  '

# Row 51
Set-TextCell $ws.Range('D51') '66.47'
$ws.Range('E51').Value = '  -0.70%  '
